# CBECC-22 / 2025 (1274) research release - initial 2025.0.1 where 2025 matches 2022
# Applies the T24R_CommunitySolar.xlsx content changes via Excel COM automation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "Last modified" mod-history note at D4 to the new entry,
#    and insert a new row (row 10) that keeps the previous note as part
#    of the Mod history list (mirrors D9's style).
# ------------------------------------------------------------------
$ws.Range("D4").Value = "10/24/22 - SAC - added 2025 code vintage records - copied directly from 2022 values for now - NEEDING UPDATE"

$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = ";"
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E9:E9").Copy() | Out-Null
$ws.Range("D10").Value = "06/01/22 - SAC - updated 2022 look-up values w/ latest 2022 kW multipliers"
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Insert a new row for the 2025 code-vintage record, right after the
#    2022 row in the CommunitySolar data table (copy formatting from the
#    2022 row, then update the year + add the "COPIED from 2022" note).
# ------------------------------------------------------------------
$ws.Rows("28:28").Insert()

$ws.Range("C27:K27").Copy()
$ws.Range("C28:K28").PasteSpecial(-4122)

$ws.Range("C28").Value = 2025
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -2137.9
$ws.Range("F28").Value = -0.0123
$ws.Range("G28").Value = -47454
$ws.Range("H28").Value = -99.150999999999996
$ws.Range("I28").Value = -1609.6
$ws.Range("J28").Value = ";"
$ws.Range("K28").Value = "SMUD Neighborhood SolarShares - Wildflower"

# Highlighted note next to the new 2025 row: bold red text on a light
# accent-4 (gold) fill, spanning P28 (text) and Q28:W28 (fill only).
$ws.Range("P28").Value = "COPIED from 2022 - MUST BE UPDATED BASED ON 2025 WEATHER & METRICS"
$ws.Range("P28:W28").Interior.ThemeColor = 8
$ws.Range("P28").Font.Bold = $true
$ws.Range("P28").Font.Color = 255
